$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the linear-model columns (E:H)
$ws.Range("E1").Value = "Coefficient"
$ws.Range("F1").Value = "Note"
$ws.Range("G1").Value = "Min_value"
$ws.Range("H1").Value = "Max_value"

# Row 2 (+++  17 .. 20)
$ws.Range("E2").Value = 0.6
$ws.Range("F2").Formula = "'+++"
$ws.Range("G2").Value = 17
$ws.Range("H2").Value = 20

# Row 3 (++  13 .. 16.5)
$ws.Range("E3").Value = 0.4
$ws.Range("F3").Formula = "'++"
$ws.Range("G3").Value = 13
$ws.Range("H3").Value = 16.5

# Row 4 (+  10 .. 12.5)
$ws.Range("F4").Value = "+"
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 12.5

# Row 5 (-  7 .. 9.5)
$ws.Range("F5").Value = "-"
$ws.Range("G5").Value = 7
$ws.Range("H5").Value = 9.5

# Row 6 (--  0 .. 6.5)
$ws.Range("F6").Formula = "'--"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 6.5

# Column widths for the new columns (target widths 15.25 / 14.75 / 14 as
# stored in the sheet1.xml <col> width; the host snaps ColumnWidth writes to
# a 1/6-character grid internally (stored = round(6*ColumnWidth)/6 + 5/6),
# so we pre-compensate by subtracting that fixed 5/6 offset to land on the
# closest achievable grid point to the target).
$ws.Columns.Item(5).ColumnWidth = 14.41666666666667
$ws.Columns.Item(6).ColumnWidth = 13.91666666666667
$ws.Columns.Item(7).ColumnWidth = 13.16666666666667

# Selection / scroll position
$ws.Range("F7").Select()
